$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 2 (the "blurred / reg" sample row) down to row 4, matching the
# existing pattern in the sheet, then tweak the single changed value.
$ws.Range("A2:L2").Copy() | Out-Null
$ws.Range("A4:L4").PasteSpecial() | Out-Null

# The only real data change: B4 becomes 65 instead of the copied 51.
$ws.Range("B4").Value = 65

# Update the active selection to match the new last-used cell.
$ws.Range("B4").Select() | Out-Null
